# Fixup of many hard coded parameters
# - corrects a couple of shared-string typos (trailing spaces / placeholder name)
# - adds two new donor rows (7 and 8) with their formatting
# - widens column Z a bit so the longer card numbers/notes are readable
# - leaves the cursor/selection on the last cell entered

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix existing shared-string values used by rows 5 and 6
# ---------------------------------------------------------------------------
$ws.Range("H5").Value = "Sam Ballard:PDQ Employee"
$ws.Range("H6").Value = "Jake Barnes"

# ---------------------------------------------------------------------------
# 2) Add new row 7
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = 42215
$ws.Range("F7").NumberFormat = "MM/DD/YY"

$ws.Range("H7").Value = "Tanya Daniels"
$ws.Range("L7").Value = "2035 Library Rd"
$ws.Range("P7").Value = "Portland, ME"
$ws.Range("T7").Value = 44236
$ws.Range("V7").Value = "(212) 555-1212"

$ws.Range("Z7").Value = 23529000445172
$ws.Range("Z7").WrapText = $true

$ws.Range("AD7").Value = "47230 - Membership Dues"
$ws.Range("AH7").Value = 40
$ws.Range("AJ7").Value = "MPL"

# ---------------------------------------------------------------------------
# 3) Add new row 8
# ---------------------------------------------------------------------------
$ws.Range("F8").Value = 42215
$ws.Range("F8").NumberFormat = "MM/DD/YY"

$ws.Range("H8").Value = "Edna Acosta"

$ws.Range("L8").Value = "7896 Library Rd. "
$ws.Range("L8").WrapText = $true

$ws.Range("P8").Value = "Portland, OR"
$ws.Range("T8").Value = 44240

$ws.Range("V8").Value = "(212) 555-1212 "
$ws.Range("V8").WrapText = $true

$ws.Range("Y8").HorizontalAlignment = -4108

$ws.Range("Z8").Value = 23529001000463
$ws.Range("Z8").WrapText = $true

$ws.Range("AA8").WrapText = $true
$ws.Range("AB8").WrapText = $true

$ws.Range("AD8").Value = "47230 - Membership Dues"
$ws.Range("AH8").Value = 10
$ws.Range("AJ8").Value = "MPL"

# ---------------------------------------------------------------------------
# 4) Widen column Z a little (account numbers got longer)
# ---------------------------------------------------------------------------
$ws.Columns.Item(26).ColumnWidth = 20.7

# ---------------------------------------------------------------------------
# 5) Leave selection on the last cell that was filled in
# ---------------------------------------------------------------------------
$ws.Range("AJ8").Select() | Out-Null
